$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: A1 = "id", B1 = "name" (C1 stays "age")
$ws.Range("A1").Value2 = "id"
$ws.Range("B1").Value2 = "name"

# Shift the A column values (x, y, z) into B column, clearing A
$ws.Range("B2").Value2 = $ws.Range("A2").Value2
$ws.Range("B3").Value2 = $ws.Range("A3").Value2
$ws.Range("B4").Value2 = $ws.Range("A4").Value2

$ws.Range("A2").Value2 = ""
$ws.Range("A3").Value2 = ""
$ws.Range("A4").Value2 = ""
